$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 8-10 (three rows drop out of the match list)
$ws.Range("A8:G10").Delete()

# Row 2
$ws.Range("A2").Value = "Full-Stack Software Engineer, Manufacturing/R&D Data Platform (NestJS, Next.js, Kafka)"
$ws.Range("B2").Value = "Sakuu Corp"
$ws.Range("C2").Value = "San Jose, CA, US USA"
$ws.Range("D2").Value = 14.4
$ws.Range("E2").Value = "FastAPI, Docker, Kubernetes, CI/CD, GitHub Actions, Git, Kafka, PostgreSQL, MongoDB, Python"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "2026-02-26"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "https://www.indeed.com/viewjob?jk=43351bc15987c6d7"

# Row 3
$ws.Range("A3").Value = "Senior Software Engineer , Backend - Dining Technology"
$ws.Range("B3").Value = "American Express"
$ws.Range("C3").Value = "New York, NY, US USA"
$ws.Range("D3").Value = 13.3
$ws.Range("E3").Value = "Docker, Kubernetes, Jenkins, GitHub Actions, Git, Kafka, MySQL, MongoDB, Python, SQL"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2026-02-20"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "https://www.indeed.com/viewjob?jk=1a55027269ae2c13"

# Row 4 (Match Score stays 12.2 - unchanged by this edit)
$ws.Range("A4").Value = "Compliance - Technology Operational Risk Management - Data Scientist - Associate"
$ws.Range("B4").Value = "JPMorganChase"
$ws.Range("C4").Value = "Plano, TX, US USA"
$ws.Range("E4").Value = "Data Scientist, RAG, TensorFlow, AWS SageMaker, Databricks, Tableau, Quicksight, Matplotlib, Python, SQL"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "2026-02-27"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "https://www.indeed.com/viewjob?jk=045a831304ccc9c5"

# Row 5
$ws.Range("A5").Value = "Senior Associate Analytics Solutions"
$ws.Range("B5").Value = "JPMorganChase"
$ws.Range("C5").Value = "New York, NY, US USA"
$ws.Range("D5").Value = 11.1
$ws.Range("E5").Value = "Generative AI, RAG, Git, Hadoop, Cassandra, Tableau, Python, SQL, R, Scala"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "2026-02-27"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "https://www.indeed.com/viewjob?jk=416f86049737d8ab"

# Row 6
$ws.Range("A6").Value = "Senior Software Engineer - AI Infrastructure"
$ws.Range("B6").Value = "Oracle"
$ws.Range("C6").Value = "Austin, TX, US USA"
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = "RAG, Docker, Terraform, NoSQL, Python, SQL, R, Java, Scala"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "2026-02-23"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = "https://www.indeed.com/viewjob?jk=6f57588aab62f26d"

# Row 7
$ws.Range("A7").Value = "Software Engineer III- Python / Numpy / Pandas"
$ws.Range("B7").Value = "JPMorganChase"
$ws.Range("C7").Value = "Jersey City, NJ, US USA"
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = "RAG, CI/CD, Jenkins, Git, Python, SQL, R, Java, Scala"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "2026-02-27"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = "https://www.indeed.com/viewjob?jk=7ba1ecac7a817c1c"
